$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.381.28"
$ws.Range("E2").Value = "  +0.02%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.569.72"
$ws.Range("E3").Value = "  +0.05%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.05%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.53"
$ws.Range("E6").Value = "  +0.46%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3762"
$ws.Range("E7").Value = "  +2.15%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.66"
$ws.Range("E8").Value = "  +0.66%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3407"
$ws.Range("E9").Value = "  +0.82%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07604"
$ws.Range("E10").Value = "  -0.22%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.142"
$ws.Range("E11").Value = "  -1.69%  "

# Row 12
$ws.Range("E12").Value = "  +0.08%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.04"
$ws.Range("E13").Value = "  -1.30%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.989"
$ws.Range("E14").Value = "  -1.12%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.965"
$ws.Range("E15").Value = "  +0.99%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.573.81"
$ws.Range("E16").Value = "  +0.03%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001132"
$ws.Range("E17").Value = "  -0.10%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.04"
$ws.Range("E18").Value = "  +0.49%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06741"
$ws.Range("E19").Value = "  -0.09%  "

# Row 20
$ws.Range("E20").Value = "  +0.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.60"
$ws.Range("E21").Value = "  +0.59%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.182"
$ws.Range("E22").Value = "  -0.60%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.95"
$ws.Range("E23").Value = "  -0.18%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.384.43"

# Row 25
$ws.Range("E25").Value = "  +0.51%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.687"
$ws.Range("E26").Value = "  -7.21%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.05"
$ws.Range("E27").Value = "  +0.25%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "147.43"
$ws.Range("E28").Value = "  +0.96%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.030"
$ws.Range("E29").Value = "  +1.45%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.40"
$ws.Range("E30").Value = "  +0.71%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.746.48"
$ws.Range("E31").Value = "  +0.10%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.016"
$ws.Range("E32").Value = "  +0.70%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9949"
$ws.Range("E33").Value = "  -4.26%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.066"
$ws.Range("E34").Value = "  -2.90%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.10"
$ws.Range("E35").Value = "  -0.64%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08450"
$ws.Range("E36").Value = "  -0.85%  "

# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.386"
$ws.Range("E37").Value = "  +10.54%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02519"
$ws.Range("E38").Value = "  -0.69%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2299"
$ws.Range("E39").Value = "  -1.00%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06500"
$ws.Range("E40").Value = "  -0.49%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.396"
$ws.Range("E41").Value = "  -2.59%  "

# Row 42
$ws.Range("E42").Value = "  -2.74%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6322"
$ws.Range("E43").Value = "  -0.53%  "

# Row 44
$ws.Range("E44").Value = "  +0.06%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.05"
$ws.Range("E45").Value = "  -2.09%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.803"
$ws.Range("E46").Value = "  +1.34%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5929"
$ws.Range("E47").Value = "  -0.86%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.278"
$ws.Range("E48").Value = "  +1.34%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.077"
$ws.Range("E49").Value = "  -1.40%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.27"
$ws.Range("E50").Value = "  -0.18%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07315"
$ws.Range("E51").Value = "  +0.33%  "
